$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers to Excel's auto-detection
# but must remain stored as text (matching the source data format).
# Force a text number format before assignment, then restore the default
# style afterwards so no stray style index is left behind.
$textCells = @("D5", "D6", "D11", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D29", "D30", "D32", "D34", "D38", "D39", "D41", "D43", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated coin data
$ws.Range("D2").Value = '66.199.50'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").Value = '3.535.35'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '608.64'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '144.18'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("D7").Value = '3.533.29'
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -4.49%  '
$ws.Range("D11").Value = '8.06'
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("E12").Value = '  -2.63%  '
$ws.Range("D13").Value = '4.133.69'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '0.0000207'
$ws.Range("E14").Value = '  -4.72%  '
$ws.Range("D15").Value = '30.44'
$ws.Range("E15").Value = '  -4.62%  '
$ws.Range("D16").Value = '3.531.18'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '66.279.44'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '10.96'
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("D20").Value = '6.24'
$ws.Range("E20").Value = '  -3.15%  '
$ws.Range("D21").Value = '15.05'
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("D22").Value = '426.03'
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("D23").Value = '0.604'
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("D24").Value = '78.88'
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("D25").Value = '3.675.09'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("E28").Value = '  -5.09%  '
$ws.Range("D29").Value = '8.05'
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("D30").Value = '2.47'
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").Value = '0.162'
$ws.Range("E32").Value = '  -3.49%  '
$ws.Range("E33").Value = '  -6.23%  '
$ws.Range("D34").Value = '25.37'
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").Value = '3.523.44'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  -2.95%  '
$ws.Range("D38").Value = '7.88'
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").Value = '5.62'
$ws.Range("E39").Value = '  -5.14%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '171.76'
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  -3.87%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '5.18'
$ws.Range("E43").Value = '  -4.86%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.894'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("E45").Value = '  -9.42%  '
$ws.Range("D46").Value = '45.23'
$ws.Range("E46").Value = '  -1.81%  '
$ws.Range("D47").Value = '26.12'
$ws.Range("E47").Value = '  -7.79%  '
$ws.Range("E48").Value = '  -6.31%  '
$ws.Range("D49").Value = '2.42'
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("D50").Value = '7.19'
$ws.Range("E50").Value = '  -4.04%  '
$ws.Range("D51").Value = '0.953'
$ws.Range("E51").Value = '  -4.00%  '

# Restore default styling on the cells we had to force to text format
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
